# Append a new paragraph at the very end of the document body, matching
# the way Word round-trips a misspelled word ("Telible") with proofErr
# (spell-check) markers around it, followed by the rest of the sentence
# in a second run.
$d = $word.ActiveDocument

$end = $d.Content.End
$r = $d.Range($end, $end)

$xml = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:rPr><w:lang w:val="es-MX"/></w:rPr></w:pPr>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t>Telible</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:lang w:val="es-MX"/></w:rPr><w:t xml:space="preserve"> lo que sucederá en el final</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$r.InsertXML($xml)
